$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.812.81'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '3.291.74'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '255.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '618.42'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.41'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +24.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.399'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.885'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.75%  '
$ws.Range('D11').Value = '3.288.49'
$ws.Range('E11').Value = '  -2.05%  '
$ws.Range('E12').Value = '  -1.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.80'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.61%  '
$ws.Range('D14').Value = '97.484.81'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000246'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.06%  '
$ws.Range('D16').Value = '3.908.66'
$ws.Range('E16').Value = '  -1.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').Value = '3.288.70'
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.50'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '475.28'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.81%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.06'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.0000202'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.54'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.51'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.74'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.66%  '
$ws.Range('D28').Value = '3.475.01'
$ws.Range('E28').Value = '  -1.69%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.288'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +19.24%  '
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.57%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '27.30'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.89%  '
$ws.Range('E36').Value = '  -3.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.09'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.56%  '
$ws.Range('E38').Value = '  -2.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.78'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '488.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.452'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.64'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.47%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.788'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.94%  '
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.12'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '159.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('E48').Value = '  -4.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.831'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.64%  '
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('E51').Value = '  -1.31%  '
